$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: update values (style/format stays as-is) ---
$ws.Range("A2").Value = 2508600
$ws.Range("B2").Value = 104989
$ws.Range("C2").Value = 12.5

# --- Rows 3 & 4: clear their contents and reset their formatting to match
#     the blank "template" rows below them (row 5's per-column styles) ---
$ws.Range("A3:E4").ClearContents()
$ws.Range("A5:E5").Copy()
$ws.Range("A3:E4").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- Rows 13 & 14: column D formatting aligned with the rest of the column
#     (copy the format already used from row 15 onward) ---
$ws.Range("D15").Copy()
$ws.Range("D13:D14").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# --- Remove the now-unused trailing blank rows 89 & 90 ---
$ws.Rows("89:90").Delete()

# --- Conditional formatting: shrink duplicate-value highlighting from
#     A2:A4 down to just A2, bumping rule priority numbers in the process
#     (mirrors Excel re-scoping the "applies to" range via the rules
#     manager, which renumbers priorities but keeps the same dxf styles) ---
$fcs = $ws.Range("A2:A4").FormatConditions
$ruleCount = $fcs.Count
for ($i = 1; $i -le $ruleCount; $i++) {
    $fcs.Item($i).ModifyAppliesToRange($ws.Range("A2"))
}
$fcs2 = $ws.Range("A2").FormatConditions
for ($i = 1; $i -le $fcs2.Count; $i++) {
    $fcs2.Item($i).Priority = 10 + $i
}

# --- Selection moves to C3 ---
$ws.Range("C3").Select()

Write-Host "edit applied"
